$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append two new rows (16, 17) at the bottom, copying row 15's formatting ---
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A17").PasteSpecial(-4122)   # xlPasteFormats

# --- The line/extr table got re-solved: two new rows (line7, line8) were
#     logically inserted after line6, shifting extr1..extr8 down. Since new
#     rows can only be appended at the bottom of the sheet, the row LABELS
#     (B) shift up into rows 8-15, and the freed bottom rows (16-17) pick up
#     the overflow (extr7, extr8). Rewrite B8:B17 and the C/D/E data to match. ---

# Row 8 -> line7 (brand-new data)
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# Row 9 -> line8 (brand-new data)
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Row 10 -> extr1 (unchanged data, relabeled)
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# Row 11 -> extr2 (unchanged data, relabeled)
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# Row 12 -> extr3 (unchanged data, relabeled)
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $true

# Row 13 -> extr4 (unchanged data, relabeled)
$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

# Row 14 -> extr5 (unchanged data, relabeled)
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

# Row 15 -> extr6 (unchanged data, relabeled)
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# Row 16 -> extr7 (new row, in_service corrected to true)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

# Row 17 -> extr8 (new row, unchanged data)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true

$wb.Save()
